$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) cells.
# D-column price cells are forced to text (NumberFormat "@") before the
# assignment so values like "206.18" or "61.36" are not silently turned
# into floating point numbers, then ClearFormats() drops the now-unneeded
# number format so the cell keeps using the sheet default style (no "s"
# attribute), matching the original inline-string cells.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.914.78'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.565.76'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.789.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.564.00'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.515'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.910.11'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.36'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '215.51'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.20'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.45%  '
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.31'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('E26').Value = '  +2.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.96'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.10'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.30%  '
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.397.73'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.925'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('E39').Value = '  +3.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.828'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.992'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  +6.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.79'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.19'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.702.21'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.90'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0980'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('E51').Value = '  +1.68%  '
